# Rename sheets:
#   "thick_thin" -> "thickthin"
#   "_x0009_острый" -> "sharp"
$wb = $excel.ActiveWorkbook

$wsThick = $wb.Sheets.Item(1)
$wsThick.Name = "thickthin"

$wsSharp = $wb.Sheets.Item(3)
$wsSharp.Name = "sharp"

# Move the active tab / selection from the 3rd sheet ("sharp") to the
# 1st sheet ("thickthin"), and update its selected cell to C18.
$wsThick.Activate()
$wsThick.Range("C18").Select()
